$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price ("D") column values are stored as text (not numbers) in this
# workbook. Force the cells to a Text number format before writing so the
# COM layer doesn't silently coerce the numeric-looking strings into
# floating point Number cells (which would also lose exact formatting like
# trailing zeros).
$priceUpdates = [ordered]@{
    "D2"  = "245.18"
    "D3"  = "23.90"
    "D4"  = "5.312"
    "D5"  = "0.05771"
    "D6"  = "6.477"
    "D7"  = "3.332"
    "D8"  = "0.8110"
    "D9"  = "0.8816"
    "D10" = "0.1391"
    "D11" = "0.07341"
    "D12" = "0.03118"
    "D13" = "0.03059"
    "D14" = "0.09344"
    "D15" = "3.870"
    "D16" = "0.001570"
    "D17" = "0.04739"
    "D18" = "0.0006014"
    "D19" = "0.005982"
    "D20" = "0.001279"
    "D22" = "0.00008813"
    "D23" = "3.604"
    "D26" = "0.1317"
    "D28" = "0.0002351"
    "D40" = "0.03804"
    "D41" = "0.006430"
    "D42" = "0.004006"
    "D43" = "0.1055"
    "D44" = "0.008369"
    "D45" = "0.00005402"
    "D47" = "0.6904"
}

foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
}

# Rows 42 and 43 swapped places (CEJI now ranks above BKEXToken), so the
# Coin name, Link and rank-prefixed Volume(1h) columns need to be updated
# to reflect the new ordering.
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("E43").Value = "42BKEXTokenBKK"
